# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Reverse the "Periodo Mora" list in column E (rows 16-23):
# old: 2212, 2301, 2302, 2303, 2304, 2305, 2306, 2307
# new: 2307, 2306, 2305, 2304, 2303, 2302, 2301, 2212
$periodos = @("2307", "2306", "2305", "2304", "2303", "2302", "2301", "2212")
for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periodos[$i]
}

# Update "Valor Mora" in column F to follow the reversed periods:
# the mora value that used to sit with period 2212 (row16 -> F=40000) now
# belongs to row 23, and the mora value that used to sit with period 2307
# (row23 -> F=16000) now belongs to row 16.
$ws.Cells.Item(16, 6).Value = 16000
$ws.Cells.Item(23, 6).Value = 40000
